$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2910.7576
$ws.Range("I15").Value = 2910.7576
$ws.Range("K15").Value = 8732.272799999999
$ws.Range("M15").Value = -8563.272799999999
$ws.Range("H98").Value = 1164.9474
$ws.Range("I98").Value = 1185.2222
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 1185.2222
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = 312.7778000000001
$ws.Range("N98").Value = -3796
$ws.Range("H112").Value = 28573284
$ws.Range("J112").Value = 2137.9644
$ws.Range("L112").Value = 6413.8932
$ws.Range("N112").Value = -8629.893199999999
$ws.Range("H122").Value = 1164.9474
$ws.Range("I122").Value = 1185.2222
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 3555.6666
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -1105.6666
$ws.Range("N122").Value = -7300
$ws.Range("H136").Value = 37113
$ws.Range("J136").Value = 37113
$ws.Range("L136").Value = 37113
$ws.Range("N136").Value = -47313
$ws.Range("H138").Value = 2285084.2
$ws.Range("I138").Value = 1467.6061
$ws.Range("J138").Value = 4169068
$ws.Range("K138").Value = 4402.8183
$ws.Range("L138").Value = 12507204
$ws.Range("M138").Value = 737.1817000000001
$ws.Range("N138").Value = -12517484

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3833005.5
$ws.Range("I122").Value = 1826.7778
$ws.Range("J122").Value = 10102207
$ws.Range("K122").Value = 5480.3334
$ws.Range("L122").Value = 30306621
$ws.Range("M122").Value = -3030.3334
$ws.Range("N122").Value = -30311521
$ws.Range("H132").Value = 299594.44
$ws.Range("I132").Value = 251040.25
$ws.Range("J132").Value = 364333.34
$ws.Range("K132").Value = 753120.75
$ws.Range("L132").Value = 1093000.02
$ws.Range("M132").Value = -750590.75
$ws.Range("N132").Value = -1098060.02

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1076.4286
$ws.Range("I20").Value = 1088
$ws.Range("K20").Value = 1088
$ws.Range("M20").Value = -841
$ws.Range("H35").Value = 20000
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20620

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7933.3335
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 7933.3335
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 7933.3335
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -8523.333500000001
$ws.Range("H34").Value = 7933.3335
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 7933.3335
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7933.3335
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -8337.333500000001
$ws.Range("H62").Value = 2421.7144
$ws.Range("I62").Value = 2089.8
$ws.Range("J62").Value = 3251.5
$ws.Range("K62").Value = 2089.8
$ws.Range("L62").Value = 3251.5
$ws.Range("M62").Value = -1465.8
$ws.Range("N62").Value = -4499.5
$ws.Range("H65").Value = 2421.7144
$ws.Range("I65").Value = 2089.8
$ws.Range("J65").Value = 3251.5
$ws.Range("K65").Value = 10449
$ws.Range("L65").Value = 16257.5
$ws.Range("M65").Value = -7329
$ws.Range("N65").Value = -22497.5
$ws.Range("H68").Value = 31168.572
$ws.Range("J68").Value = 33863.332
$ws.Range("L68").Value = 33863.332
$ws.Range("N68").Value = -35361.332
$ws.Range("H71").Value = 31168.572
$ws.Range("J71").Value = 33863.332
$ws.Range("L71").Value = 101589.996
$ws.Range("N71").Value = -109077.996
$ws.Range("H99").Value = 9145.6
$ws.Range("I99").Value = 9098.666999999999
$ws.Range("J99").Value = 9333.333000000001
$ws.Range("K99").Value = 9098.666999999999
$ws.Range("L99").Value = 9333.333000000001
$ws.Range("M99").Value = -7600.666999999999
$ws.Range("N99").Value = -12329.333
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H126").Value = 9145.6
$ws.Range("I126").Value = 9098.666999999999
$ws.Range("J126").Value = 9333.333000000001
$ws.Range("K126").Value = 27296.001
$ws.Range("L126").Value = 27999.999
$ws.Range("M126").Value = -24826.001
$ws.Range("N126").Value = -32939.999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 479.34616
$ws.Range("I5").Value = 366.4737
$ws.Range("J5").Value = 785.7143
$ws.Range("K5").Value = 1099.4211
$ws.Range("L5").Value = 2357.1429
$ws.Range("M5").Value = -987.4211
$ws.Range("N5").Value = -2581.1429
$ws.Range("H11").Value = 559093.5
$ws.Range("I11").Value = 375382.5
$ws.Range("J11").Value = 1000000
$ws.Range("K11").Value = 1126147.5
$ws.Range("L11").Value = 3000000
$ws.Range("M11").Value = -1126007.5
$ws.Range("N11").Value = -3000280
$ws.Range("H122").Value = 596.53845
$ws.Range("I122").Value = 267.5
$ws.Range("K122").Value = 2407.5
$ws.Range("M122").Value = 42.5
$ws.Range("H131").Value = 1171
$ws.Range("J131").Value = 1251.6857
$ws.Range("L131").Value = 3755.0571
$ws.Range("N131").Value = -13835.0571
$ws.Range("H135").Value = 479.34616
$ws.Range("I135").Value = 366.4737
$ws.Range("J135").Value = 785.7143
$ws.Range("K135").Value = 3298.2633
$ws.Range("L135").Value = 7071.428699999999
$ws.Range("M135").Value = -763.2633000000001
$ws.Range("N135").Value = -12141.4287

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 88.916664
$ws.Range("I2").Value = 98.59999999999999
$ws.Range("J2").Value = 82
$ws.Range("K2").Value = 98.59999999999999
$ws.Range("L2").Value = 82
$ws.Range("M2").Value = 14.40000000000001
$ws.Range("N2").Value = -308
$ws.Range("H10").Value = 20000000
$ws.Range("I10").Value = 20000000
$ws.Range("K10").Value = 20000000
$ws.Range("M10").Value = -19999831
$ws.Range("H70").Value = 42229.15
$ws.Range("I70").Value = 67242.94
$ws.Range("J70").Value = 5845.4546
$ws.Range("K70").Value = 67242.94
$ws.Range("L70").Value = 5845.4546
$ws.Range("M70").Value = -66972.94
$ws.Range("N70").Value = -6385.4546
$ws.Range("H73").Value = 42229.15
$ws.Range("I73").Value = 67242.94
$ws.Range("J73").Value = 5845.4546
$ws.Range("K73").Value = 67242.94
$ws.Range("L73").Value = 5845.4546
$ws.Range("M73").Value = -66306.94
$ws.Range("N73").Value = -7717.4546
$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -6400
$ws.Range("H123").Value = 26078.059
$ws.Range("J123").Value = 26078.059
$ws.Range("L123").Value = 26078.059
$ws.Range("N123").Value = -30978.059
$ws.Range("H132").Value = 62187.668
$ws.Range("I132").Value = 37221.18
$ws.Range("J132").Value = 202000
$ws.Range("K132").Value = 111663.54
$ws.Range("L132").Value = 606000
$ws.Range("M132").Value = -109133.54
$ws.Range("N132").Value = -611060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2193.6
$ws.Range("I7").Value = 2174.0833
$ws.Range("J7").Value = 2222.875
$ws.Range("K7").Value = 2174.0833
$ws.Range("L7").Value = 2222.875
$ws.Range("M7").Value = -2062.0833
$ws.Range("N7").Value = -2446.875
$ws.Range("H40").Value = 3458.5833
$ws.Range("I40").Value = 3333.8333
$ws.Range("J40").Value = 3583.3333
$ws.Range("K40").Value = 3333.8333
$ws.Range("L40").Value = 3583.3333
$ws.Range("M40").Value = -3197.8333
$ws.Range("N40").Value = -3855.3333
$ws.Range("H55").Value = 120.6875
$ws.Range("J55").Value = 170.4
$ws.Range("L55").Value = 170.4
$ws.Range("N55").Value = -516.4
$ws.Range("H122").Value = 2853.7104
$ws.Range("I122").Value = 2211.625
$ws.Range("J122").Value = 3320.682
$ws.Range("K122").Value = 6634.875
$ws.Range("L122").Value = 9962.045999999998
$ws.Range("M122").Value = -4184.875
$ws.Range("N122").Value = -14862.046
$ws.Range("H126").Value = 2193.6
$ws.Range("I126").Value = 2174.0833
$ws.Range("J126").Value = 2222.875
$ws.Range("K126").Value = 6522.249899999999
$ws.Range("L126").Value = 6668.625
$ws.Range("M126").Value = -4052.249899999999
$ws.Range("N126").Value = -11608.625
$ws.Range("H136").Value = 72738.664
$ws.Range("I136").Value = 69927.19
$ws.Range("J136").Value = 75951.78999999999
$ws.Range("K136").Value = 209781.57
$ws.Range("L136").Value = 227855.37
$ws.Range("M136").Value = -207231.57
$ws.Range("N136").Value = -232955.37

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 30706.307
$ws.Range("J76").Value = 30706.307
$ws.Range("L76").Value = 30706.307
$ws.Range("N76").Value = -31336.307
$ws.Range("H79").Value = 30706.307
$ws.Range("J79").Value = 30706.307
$ws.Range("L79").Value = 30706.307
$ws.Range("N79").Value = -32890.307
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H122").Value = 2996.64
$ws.Range("I122").Value = 2601.1333
$ws.Range("J122").Value = 3589.9
$ws.Range("K122").Value = 7803.3999
$ws.Range("L122").Value = 10769.7
$ws.Range("M122").Value = -5353.3999
$ws.Range("N122").Value = -15669.7
$ws.Range("H126").Value = 1272.7059
$ws.Range("I126").Value = 1052
$ws.Range("J126").Value = 1990
$ws.Range("K126").Value = 3156
$ws.Range("L126").Value = 5970
$ws.Range("M126").Value = -686
$ws.Range("N126").Value = -10910
$ws.Range("H132").Value = 37037.824
$ws.Range("I132").Value = 27909.59
$ws.Range("J132").Value = 92567.914
$ws.Range("K132").Value = 83728.77
$ws.Range("L132").Value = 277703.742
$ws.Range("M132").Value = -81198.77
$ws.Range("N132").Value = -282763.742
$ws.Range("H136").Value = 120099.06
$ws.Range("I136").Value = 125960.5
$ws.Range("J136").Value = 114888.89
$ws.Range("K136").Value = 377881.5
$ws.Range("L136").Value = 344666.67
$ws.Range("M136").Value = -375331.5
$ws.Range("N136").Value = -349766.67
$ws.Range("H140").Value = 49342.125
$ws.Range("J140").Value = 49342.125
$ws.Range("L140").Value = 49342.125
$ws.Range("N140").Value = -59702.125
$ws.Range("H141").Value = 57969.7
$ws.Range("J141").Value = 57969.7
$ws.Range("L141").Value = 57969.7
$ws.Range("N141").Value = -68329.7
